$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168, shifting existing rows 168..291 down to 169..292
$ws.Rows(168).Insert()

# Populate the newly inserted row 168 with the new Jengibre price record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R mirror the constant values used throughout
# this product block; D,J,K,L,M,P carry the row-specific figures.
$ws.Cells.Item(168, 1).Value = 10
$ws.Cells.Item(168, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(168, 3).Value = "La Araucanía"
$ws.Cells.Item(168, 4).Value = 45068
$ws.Cells.Item(168, 5).Value = 9
$ws.Cells.Item(168, 6).Value = 100114007
$ws.Cells.Item(168, 7).Value = "Jengibre"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 100
$ws.Cells.Item(168, 11).Value = 24000
$ws.Cells.Item(168, 12).Value = 24000
$ws.Cells.Item(168, 13).Value = 24000
$ws.Cells.Item(168, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(168, 15).Value = "Perú"
$ws.Cells.Item(168, 16).Value = 1846
$ws.Cells.Item(168, 17).Value = 13
$ws.Cells.Item(168, 18).Value = "Hortaliza"
